# Work plan.xlsx -- add detection / BER calculation results
# Mirrors the authored commit: mark "Coherent detection" & "Binary detector"
# as Done, record the BER-calculation subtask, and change the informal
# "Actual Time" decimals (0.5 / 3) into readable text (30 min / 3 h / 1 h).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dashboard")

# --- New text entered first (drives the shared-string append order) ---
$ws.Range("H10").Value = "Need BER calculation "
$ws.Range("G9").Value  = "1 h"
$ws.Range("G6").Value  = "3 h"
$ws.Range("G2").Value  = "30 min"
$ws.Range("G3").Value  = "30 min"
$ws.Range("G4").Value  = "30 min"
$ws.Range("G5").Value  = "30 min"
$ws.Range("G10").Value = "30 min"

# --- Finish filling in the "Binary detector" row (BER calculation task) ---
$ws.Range("F10").Value = "23/09/21"

# --- Status updates: both tasks are now complete ---
$ws.Range("I9").Value  = "Done"
$ws.Range("I10").Value = "Done"

# --- Row 3 ("Create a project and automate ...") gets its cells vertically
#     centered, matching the taller (ht=30) wrapped row above/below it ---
$ws.Range("B3:J3").VerticalAlignment = -4108

# --- Leave the selection where the author left it ---
$ws.Range("H2").Select()
